$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.401.38'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.572.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.35'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3767'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.90'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.158'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07650'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.26'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.949'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.571.98'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001133'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.24'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06760'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.82'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.209'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.428'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '22.394.61'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.699'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.26'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.039'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.747.60'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.167'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.008'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9965'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.986'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08553'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02549'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2316'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06585'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.325'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.453'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.51'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.10%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6422'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.01'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.793'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5999'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.304'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +7.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.087'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.66'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("E51").Value = '  +0.52%  '
